$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# 1. Update I4 on Typography sheet: "0-9,a-z" -> "0-9,a-z,A-Z"
$wsTypo.Range("I4").Value = "0-9,a-z,A-Z"

# 2. Fill in row 23 on Translation sheet
$wsTrans.Range("B23").Value = "SingleUseId22"
$wsTrans.Range("C23").Value = "Default"
$wsTrans.Range("D23").Value = "Center"
$wsTrans.Range("E23").Value = "LTR"
$wsTrans.Range("F23").Value = "<text>"
